$d = $word.ActiveDocument

# 1. Title
$d.Content.Find.Execute(
    "Beyond Sight: The Realm of Ultraviolet Perception",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The Intriguing World of Chemistry: An Exploration of Matter and Its Properties",
    2) | Out-Null

# 2. Author name: "Dr. Eliana Conti" -> "Richard Henderson"
$d.Content.Find.Execute(
    "Dr. Eliana Conti",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Richard Henderson",
    2) | Out-Null

# 3. Email: "eliana.conti@astroscience.org" -> "richardhenderson@gmail.net"
$d.Content.Find.Execute(
    "eliana.conti@astroscience.org",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "richardhenderson@gmail.net",
    2) | Out-Null

# 4. Body paragraph, segment A (before the first double line-break)
$segA_old = "In the boundless tapestry of the natural world, humans perceive but a fragmented glimpse of the vibrant symphony of light that weaves its way through existence. Our limited vision, confined within the narrow spectrum of visible light, leaves us oblivious to the vast realms that lie hidden beyond our perception. Among these unseen realms resides the enigmatic domain of ultraviolet (UV) radiation, a captivating realm where untold secrets beckon and awe-inspiring phenomena unfold. In this essay, we will embark on a journey into the fascinating world of UV perception, unveiling its profound implications for diverse fields of study, from astronomy and biology to engineering and medicine, revealing how this hidden realm expands our understanding of the universe and ourselves."
$segA_new = "Chemistry, the study of matter and its composition, properties, and change, offers a captivating exploration of the universe around us. From the vibrant colors of nature to the intricate processes occurring within our bodies, chemistry plays a vital role in defining our everyday existence. In this essay, we delve into the fascinating realm of chemistry, unraveling the secrets of matter and discovering the astounding phenomena that govern its behavior."
$d.Content.Find.Execute(
    $segA_old,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $segA_new,
    2) | Out-Null

# 5. Body paragraph, segment B (between the two double line-breaks)
$segB_old = "As we venture into the realm of UV perception, we encounter an array of intriguing applications that span the boundaries of various disciplines. From uncovering the hidden patterns and structures of ancient artifacts to detecting counterfeit works of art, UV radiation serves as an invaluable tool for scientists and art historians alike. Likewise, in the vast expanse of astronomy, UV observations have unveiled hidden galaxies, quasars, and nebulae pulsating with celestial energy, allowing us to gain unprecedented insights into the evolution of the cosmos. In the realm of biology, UV perception has illuminated intricate communication strategies employed by various species, such as the beguiling courtship dances of fireflies, the subtle markings on butterfly wings, and the patterns on the shells of certain marine creatures. Through UV perception, we unlock a secret language of nature, where visual cues and signals shape the intricate tapestry of life."
$segB_new = "As we embark on this journey, we will uncover the basic building blocks of matter, exploring the structure of atoms and molecules. We will witness the transformative power of chemical reactions, witnessing how substances rearrange themselves to form new compounds with unique properties. Moreover, we will investigate the interactions between matter and energy, delving into the intricacies of chemical bonding and the energy changes that accompany chemical processes."
$d.Content.Find.Execute(
    $segB_old,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $segB_new,
    2) | Out-Null

# 6. Body paragraph, segment C (after the second double line-break)
$segC_old = "Beyond its applications in scientific research, UV perception has also found its niche in everyday technologies and medical advancements. UV lamps find widespread use in sterilization and disinfection, effectively combating harmful bacteria and viruses in hospitals, clinics, and food processing facilities. In the realm of engineering, UV radiation plays a vital role in non-destructive testing, revealing hidden cracks and defects in materials and structures. Furthermore, the advent of UV-based treatments has revolutionized various medical specialties, such as dermatology, where UV radiation is harnessed to treat skin conditions like psoriasis and vitiligo. In ophthalmology, UV-A light therapy offers promising solutions for managing ocular ailments. As we delve deeper into the captivating universe of UV perception, its profound implications continue to inspire and intrigue, pushing the boundaries of our knowledge and transforming our understanding of the world around us."
$segC_new = "Beyond the fundamental aspects of chemistry, we will venture into the practical applications that have shaped our world. From the creation of new materials to the development of life-saving medicines, chemistry has revolutionized numerous industries and improved countless lives. We will examine the role of chemistry in our daily lives, examining the chemistry behind common household products and unraveling the mysteries of biological processes."
$d.Content.Find.Execute(
    $segC_old,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $segC_new,
    2) | Out-Null

# 7. Summary paragraph
$sumA_old = "In the realm of UV perception, we have explored the myriad ways in which this hidden realm has transformed our understanding of the cosmos, biology, technology, and medicine."
$sumA_new = "Through our exploration of chemistry, we have gained a deeper understanding of the matter that constitutes our universe, the transformations it undergoes, and the profound impact it has on our lives."
$d.Content.Find.Execute(
    $sumA_old,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $sumA_new,
    2) | Out-Null

$sumB_old = " UV radiation has empowered us to unveil celestial wonders, decode the secrets of nature's communication, enhance everyday technologies, and revolutionize medical treatments."
$sumB_new = " We have witnessed the power of chemical reactions, marveled at the intricacies of atomic structure, and uncovered the practical applications that have revolutionized our world."
$d.Content.Find.Execute(
    $sumB_old,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $sumB_new,
    2) | Out-Null

$sumC_old = " As we continue to unravel the mysteries of this unseen realm, we unlock the potential for further advancements, fostering a future where the boundaries of science, art, and technology seamlessly intertwine."
$sumC_new = " This journey into the realm of chemistry serves as a testament to the boundless curiosity and ingenuity of humankind, reminding us that the pursuit of knowledge holds the key to unlocking the secrets of our physical world."
$d.Content.Find.Execute(
    $sumC_old,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $sumC_new,
    2) | Out-Null

# 8. Add a new empty paragraph at the very end of the document (after the Summary paragraph)
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Host "Done"
